$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data row 3 (the placeholder "100" row), shrinking the
# used range back down to A1:H2.
$ws.Rows.Item(3).Delete() | Out-Null

# Replace row 2's placeholder header-like text with real bird record data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Golden"
$ws.Range("C2").Value = "Amrican"

$ws.Range("D2").NumberFormat = "mm-dd-yy"
$dob = Get-Date -Year 1997 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("D2").Value = $dob.Date

$ws.Range("E2").Value = "Female"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 2
